$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_val data (filtering save games) - new literal values for B2:G8
$data = @(
    @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 4.429675500412797),
    @(0.6753301551942219, 0.3127903958511391, 3.900430680208489, 0.496779210170732, 0, 5.385330441424582),
    @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 6.201049113329182),
    @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 3.645393585217082),
    @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 1, 14.36450238910742),
    @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 8.660232485948974, 0, 11.945164432584),
    @(3.230985683306322, 1.667794583268128, 26.21740644021617, 0.496779210170732, 0, 31.61296591696135)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
